# #3495 - PurchaseInvoicePlan: add Quantity/Unit Price columns (after "Fiscal year
# by Invoice Plan") and WA Number/Recieive Quantity/Unit Price/Subtotal columns
# (before "Billing Number") to the Purchase Invoice Plan Report template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Insert two new columns at Y:Z -> "Quantity", "Unit Price"
#    (shifts the old Installment..FinLease block from Y.. to AA..)
# -----------------------------------------------------------------
$ws.Columns("Y:Z").Insert()
$ws.Range("Y10").Value = "Quantity"
$ws.Range("Z10").Value = "Unit Price"

# -----------------------------------------------------------------
# 2) Insert four new columns right before "Billing Number" (now at AN)
#    -> "WA Number", "Recieive Quantity", "Unit Price", "Subtotal"
# -----------------------------------------------------------------
$ws.Columns("AN:AQ").Insert()
$ws.Range("AN10").Value = "WA Number"
$ws.Range("AO10").Value = "Recieive Quantity"
$ws.Range("AP10").Value = "Unit Price"
$ws.Range("AQ10").Value = "Subtotal"

# -----------------------------------------------------------------
# 3) Column width / best-fit touch-ups (best effort match of the
#    resized columns in the final template).
# -----------------------------------------------------------------
$ws.Range("E1").EntireColumn.ColumnWidth = 4.857142857142857

$ws.Range("Y1").EntireColumn.ColumnWidth = 18.857142857142858
$ws.Range("Z1").EntireColumn.ColumnWidth = 18.857142857142858
$ws.Range("AA1").EntireColumn.ColumnWidth = 11.142857142857142
$ws.Range("AC1").EntireColumn.ColumnWidth = 15.571428571428571
$ws.Range("AD1").EntireColumn.ColumnWidth = 16.285714285714285
$ws.Range("AE1").EntireColumn.ColumnWidth = 16.857142857142858
$ws.Range("AI1").EntireColumn.ColumnWidth = 16.428571428571427
$ws.Range("AJ1").EntireColumn.ColumnWidth = 14.285714285714286
$ws.Range("AK1").EntireColumn.ColumnWidth = 17.142857142857142
$ws.Range("AL1").EntireColumn.ColumnWidth = 18.285714285714285
$ws.Range("AN1").EntireColumn.ColumnWidth = 15.714285714285714
$ws.Range("AO1").EntireColumn.ColumnWidth = 17.428571428571427
$ws.Range("AP1").EntireColumn.ColumnWidth = 16.857142857142858
$ws.Range("AQ1").EntireColumn.ColumnWidth = 17.714285714285715
$ws.Range("AR1").EntireColumn.ColumnWidth = 17.714285714285715
$ws.Range("AS1").EntireColumn.ColumnWidth = 17.142857142857142
$ws.Range("AU1").EntireColumn.ColumnWidth = 12.285714285714286
$ws.Range("AV1").EntireColumn.ColumnWidth = 24.0
$ws.Range("AW1").EntireColumn.ColumnWidth = 21.714285714285715
$ws.Range("AX1").EntireColumn.ColumnWidth = 20.857142857142858
$ws.Range("AY1").EntireColumn.ColumnWidth = 7.571428571428571
$ws.Range("AZ1").EntireColumn.ColumnWidth = 16.857142857142858
$ws.Range("BA1").EntireColumn.ColumnWidth = 8.857142857142858

# -----------------------------------------------------------------
# 4) Restore the cursor/selection the author left the sheet on.
# -----------------------------------------------------------------
$ws.Range("AR16").Select()
